$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the monthly values for hansraj (row 2), keep only the name
$ws.Range("B2:H2").ClearContents()

# Clear the stray value in B3 / B4 for riguda and chikne
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
